$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","Q","R")

# --- Swap rows 7 and 8 ---
$row7 = @{}
$row8 = @{}
foreach ($c in $cols) {
    $row7[$c] = $ws.Range("$c" + "7").Value()
    $row8[$c] = $ws.Range("$c" + "8").Value()
}
foreach ($c in $cols) {
    $ws.Range("$c" + "7").Value = $row8[$c]
    $ws.Range("$c" + "8").Value = $row7[$c]
}

# --- Rotate rows 21, 22, 23, 24 ---
# New row21 = old row22, new row22 = old row23, new row23 = old row24, new row24 = old row21
$row21 = @{}
$row22 = @{}
$row23 = @{}
$row24 = @{}
foreach ($c in $cols) {
    $row21[$c] = $ws.Range("$c" + "21").Value()
    $row22[$c] = $ws.Range("$c" + "22").Value()
    $row23[$c] = $ws.Range("$c" + "23").Value()
    $row24[$c] = $ws.Range("$c" + "24").Value()
}
foreach ($c in $cols) {
    $ws.Range("$c" + "21").Value = $row22[$c]
    $ws.Range("$c" + "22").Value = $row23[$c]
    $ws.Range("$c" + "23").Value = $row24[$c]
    $ws.Range("$c" + "24").Value = $row21[$c]
}
